# Apply "repull data" style updates to column F (dSF) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new F-column value
$updates = @{
    2  = -10
    11 = -4
    17 = 3
    18 = -3
    21 = -5
    27 = -1
    28 = -3
    29 = 0
    37 = -3
    39 = -2
    41 = 1
    42 = -1
    53 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
